$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1627906976744186
$ws.Range("C2").Value = 0.6186046511627907
$ws.Range("J2").Value = 0.0186046511627907
$ws.Range("O2").Value = 0.004651162790697674
$ws.Range("P2").Value = 0.1441860465116279
$ws.Range("S2").Value = 0.05116279069767442

# Row 3
$ws.Range("B3").Value = 0.007194244604316547
$ws.Range("C3").Value = 0.02877697841726619
$ws.Range("J3").Value = 0.05035971223021583
$ws.Range("P3").Value = 0.7482014388489209
$ws.Range("S3").Value = 0.1654676258992806

# Row 4
$ws.Range("J4").Value = 0.02173913043478261
$ws.Range("O4").Value = 0.02173913043478261
$ws.Range("P4").Value = 0.6956521739130435
$ws.Range("S4").Value = 0.2608695652173913

# Row 6
$ws.Range("B6").Value = 0.02008032128514056
$ws.Range("D6").Value = 0.02409638554216868
$ws.Range("F6").Value = 0.06425702811244979
$ws.Range("J6").Value = 0.2369477911646586
$ws.Range("O6").Value = 0.0321285140562249
$ws.Range("Q6").Value = 0.1967871485943775
$ws.Range("R6").Value = 0.08032128514056225
$ws.Range("S6").Value = 0.3453815261044177

# Row 7
$ws.Range("B7").Value = 0.09268292682926829
$ws.Range("D7").Value = 0.00975609756097561
$ws.Range("F7").Value = 0.05853658536585366
$ws.Range("J7").Value = 0.1268292682926829
$ws.Range("O7").Value = 0.04878048780487805
$ws.Range("Q7").Value = 0.2
$ws.Range("R7").Value = 0.08780487804878048
$ws.Range("S7").Value = 0.375609756097561

# Row 8
$ws.Range("B8").Value = 0.08967391304347826
$ws.Range("D8").Value = 0.02173913043478261
$ws.Range("F8").Value = 0.05706521739130434
$ws.Range("J8").Value = 0.1141304347826087
$ws.Range("O8").Value = 0.01630434782608696
$ws.Range("Q8").Value = 0.2065217391304348
$ws.Range("R8").Value = 0.09782608695652174
$ws.Range("S8").Value = 0.3967391304347826

# Row 9
$ws.Range("B9").Value = 0.1
$ws.Range("D9").Value = 0.01176470588235294
$ws.Range("E9").Value = 0.005882352941176471
$ws.Range("F9").Value = 0.07647058823529412
$ws.Range("J9").Value = 0.08235294117647059
$ws.Range("O9").Value = 0.03529411764705882
$ws.Range("Q9").Value = 0.2294117647058823
$ws.Range("R9").Value = 0.1294117647058824
$ws.Range("S9").Value = 0.3294117647058823

# Row 10
$ws.Range("B10").Value = 0.07959022852639874
$ws.Range("D10").Value = 0.02285263987391647
$ws.Range("F10").Value = 0.08274231678486997
$ws.Range("J10").Value = 0.1142631993695823
$ws.Range("O10").Value = 0.01812450748620961
$ws.Range("Q10").Value = 0.2379826635145784
$ws.Range("R10").Value = 0.1000788022064618
$ws.Range("S10").Value = 0.3443656422379827

# Row 11
$ws.Range("G11").Value = 0.145631067961165
$ws.Range("J11").Value = 0.09061488673139159
$ws.Range("K11").Value = 0.2103559870550162
$ws.Range("L11").Value = 0.517799352750809
$ws.Range("S11").Value = 0.03559870550161812

# Row 12
$ws.Range("G12").Value = 0.7469879518072289
$ws.Range("J12").Value = 0.2048192771084337
$ws.Range("K12").Value = 0.01807228915662651
$ws.Range("L12").Value = 0.02409638554216868
$ws.Range("S12").Value = 0.006024096385542169

# Row 13
$ws.Range("G13").Value = 0.7169811320754716
$ws.Range("J13").Value = 0.2830188679245283

# Row 15
$ws.Range("F15").Value = 0.02641509433962264
$ws.Range("H15").Value = 0.1471698113207547
$ws.Range("I15").Value = 0.07169811320754717
$ws.Range("J15").Value = 0.3773584905660378
$ws.Range("K15").Value = 0.06415094339622641
$ws.Range("M15").Value = 0.01509433962264151
$ws.Range("O15").Value = 0.0830188679245283
$ws.Range("S15").Value = 0.2150943396226415

# Row 16
$ws.Range("F16").Value = 0.03184713375796178
$ws.Range("H16").Value = 0.1401273885350318
$ws.Range("I16").Value = 0.06369426751592357
$ws.Range("J16").Value = 0.4267515923566879
$ws.Range("K16").Value = 0.09554140127388536
$ws.Range("M16").Value = 0.03184713375796178
$ws.Range("O16").Value = 0.07006369426751592
$ws.Range("S16").Value = 0.1401273885350318

# Row 17
$ws.Range("F17").Value = 0.02385685884691849
$ws.Range("H17").Value = 0.1332007952286282
$ws.Range("I17").Value = 0.08548707753479125
$ws.Range("J17").Value = 0.4294234592445328
$ws.Range("K17").Value = 0.1053677932405567
$ws.Range("M17").Value = 0.02584493041749503
$ws.Range("N17").Value = 0.003976143141153081
$ws.Range("O17").Value = 0.09940357852882704
$ws.Range("S17").Value = 0.09343936381709742

# Row 18
$ws.Range("F18").Value = 0.01339285714285714
$ws.Range("H18").Value = 0.1383928571428572
$ws.Range("I18").Value = 0.09375
$ws.Range("J18").Value = 0.4732142857142857
$ws.Range("K18").Value = 0.07142857142857142
$ws.Range("M18").Value = 0.03125
$ws.Range("O18").Value = 0.08035714285714286
$ws.Range("S18").Value = 0.09821428571428571

# Row 19
$ws.Range("F19").Value = 0.0261437908496732
$ws.Range("H19").Value = 0.1979458450046685
$ws.Range("I19").Value = 0.07096171802054155
$ws.Range("J19").Value = 0.3940242763772175
$ws.Range("K19").Value = 0.1297852474323063
$ws.Range("M19").Value = 0.02334267040149393
$ws.Range("N19").Value = 0.002801120448179272
$ws.Range("O19").Value = 0.07936507936507936
$ws.Range("S19").Value = 0.07563025210084033
